$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177, shifting rows 177:280 down to 178:281
$ws.Rows.Item(177).Insert()

# Fill in the newly inserted row 177 with the new data
$ws.Cells.Item(177, 1).Value = 7
$ws.Cells.Item(177, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(177, 3).Value = "Ñuble"
$ws.Cells.Item(177, 4).Value = 44957
$ws.Cells.Item(177, 5).Value = 16
$ws.Cells.Item(177, 6).Value = 100112032
$ws.Cells.Item(177, 7).Value = "Zapallo italiano"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 80
$ws.Cells.Item(177, 11).Value = 5000
$ws.Cells.Item(177, 12).Value = 5500
$ws.Cells.Item(177, 13).Value = 5250
$ws.Cells.Item(177, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(177, 15).Value = "Región del Maule"
$ws.Cells.Item(177, 16).Value = 105
$ws.Cells.Item(177, 17).Value = 50
$ws.Cells.Item(177, 18).Value = "Hortaliza"
